$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NIK and Nama Karyawan values (shared across C2:C7 and D2:D7)
$ws.Range("C2:C7").Value = "EN-4-047"
$ws.Range("D2:D7").Value = "Ari Pratama"

# Update the active selection shown when the sheet is opened
$ws.Range("C3:D7").Select()
